$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

# Row 377
Set-TextValue $ws.Range("B377") "7309340"
$ws.Range("E377").Value = "Lecco"
$ws.Range("F377").Value = "Modena"
$ws.Range("J377").Value = 3.8
$ws.Range("K377").Value = 3.5
$ws.Range("L377").Value = 1.95
$ws.Range("M377").Value = 3.6
$ws.Range("N377").Value = 3.5
$ws.Range("O377").Value = 2
$ws.Range("P377").Value = 0.5
$ws.Range("Q377").Value = 1.825
$ws.Range("R377").Value = 2.025
$ws.Range("S377").Value = 2.5
$ws.Range("T377").Value = 1.875
$ws.Range("U377").Value = 1.975

# Row 378
Set-TextValue $ws.Range("B378") "7309341"
$ws.Range("E378").Value = "AC Reggiana"
$ws.Range("F378").Value = "Parma"
$ws.Range("J378").Value = 3.3
$ws.Range("K378").Value = 3.3
$ws.Range("L378").Value = 2.15
$ws.Range("M378").Value = 3.4
$ws.Range("N378").Value = 3.1
$ws.Range("O378").Value = 2.25
$ws.Range("P378").Value = 0.25
$ws.Range("Q378").Value = 1.95
$ws.Range("R378").Value = 1.9
$ws.Range("S378").Value = 2.5
$ws.Range("T378").Value = 1.925
$ws.Range("U378").Value = 1.925

# Row 379
$ws.Range("M379").Value = 2.75
$ws.Range("O379").Value = 2.5
$ws.Range("Q379").Value = 2
$ws.Range("R379").Value = 1.85

# Row 380
Set-TextValue $ws.Range("B380") "7310531"
$ws.Range("E380").Value = "Cremonese"
$ws.Range("F380").Value = "Cittadella"
$ws.Range("J380").Value = 1.666
$ws.Range("K380").Value = 4
$ws.Range("L380").Value = 4.75
$ws.Range("M380").Value = 1.666
$ws.Range("N380").Value = 4
$ws.Range("O380").Value = 4.5
$ws.Range("P380").Value = -0.75
$ws.Range("Q380").Value = 1.875
$ws.Range("R380").Value = 1.975
$ws.Range("S380").Value = 2.75
$ws.Range("T380").Value = 1.925
$ws.Range("U380").Value = 1.925

# Row 381
Set-TextValue $ws.Range("B381") "7286345"
$ws.Range("E381").Value = "Bari"
$ws.Range("F381").Value = "Brescia"
$ws.Range("J381").Value = 1.6
$ws.Range("K381").Value = 3.8
$ws.Range("L381").Value = 5.75
$ws.Range("M381").Value = 1.571
$ws.Range("N381").Value = 4
$ws.Range("O381").Value = 5.5
$ws.Range("P381").Value = -1
$ws.Range("T381").Value = 2.025
$ws.Range("U381").Value = 1.825

# Row 382
$ws.Range("M382").Value = 2.875
$ws.Range("N382").Value = 3.5
$ws.Range("O382").Value = 2.3
$ws.Range("Q382").Value = 1.8
$ws.Range("R382").Value = 2.05

# Row 383
Set-TextValue $ws.Range("B383") "7280956"
$ws.Range("E383").Value = "Spezia"
$ws.Range("F383").Value = "Venezia"
$ws.Range("J383").Value = 2.25
$ws.Range("K383").Value = 3.4
$ws.Range("L383").Value = 3.1
$ws.Range("M383").Value = 2.375
$ws.Range("N383").Value = 3.4
$ws.Range("O383").Value = 2.9
$ws.Range("P383").Value = -0.25
$ws.Range("Q383").Value = 2.025
$ws.Range("R383").Value = 1.825
$ws.Range("S383").Value = 2.5
$ws.Range("T383").Value = 1.85
$ws.Range("U383").Value = 2

# Row 384
Set-TextValue $ws.Range("B384") "7280955"
$ws.Range("E384").Value = "AC Feralpisalo"
$ws.Range("F384").Value = "Ternana"
$ws.Range("J384").Value = 4
$ws.Range("K384").Value = 3.8
$ws.Range("L384").Value = 1.8
$ws.Range("M384").Value = 4.5
$ws.Range("N384").Value = 4
$ws.Range("O384").Value = 1.666
$ws.Range("P384").Value = 0.75
$ws.Range("Q384").Value = 1.975
$ws.Range("R384").Value = 1.875
$ws.Range("S384").Value = 3
$ws.Range("T384").Value = 1.975
$ws.Range("U384").Value = 1.875

# Row 385
Set-TextValue $ws.Range("B385") "7280954"
$ws.Range("E385").Value = "Como"
$ws.Range("F385").Value = "Cosenza"
$ws.Range("J385").Value = 1.363
$ws.Range("K385").Value = 4.75
$ws.Range("L385").Value = 8
$ws.Range("M385").Value = 1.333
$ws.Range("N385").Value = 5
$ws.Range("O385").Value = 9
$ws.Range("P385").Value = -1.5
$ws.Range("Q385").Value = 2.025
$ws.Range("R385").Value = 1.825
$ws.Range("S385").Value = 2.75
$ws.Range("T385").Value = 1.825
$ws.Range("U385").Value = 2.025

# Row 386
Set-TextValue $ws.Range("B386") "7280953"
$ws.Range("E386").Value = "Ascoli"
$ws.Range("F386").Value = "Pisa"
$ws.Range("J386").Value = 1.666
$ws.Range("K386").Value = 3.75
$ws.Range("L386").Value = 5
$ws.Range("M386").Value = 1.65
$ws.Range("N386").Value = 3.75
$ws.Range("O386").Value = 5.25
$ws.Range("P386").Value = -0.75
$ws.Range("Q386").Value = 1.85
$ws.Range("R386").Value = 2
$ws.Range("S386").Value = 2.25
$ws.Range("T386").Value = 2.025
$ws.Range("U386").Value = 1.825
